$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Fecha" (D) and "Volumen" (M) columns per the weekly re-sequencing.
# Row 3
$ws.Range("D3").Value = 44313
$ws.Range("M3").Value = 120

# Row 4
$ws.Range("D4").Value = 44302
$ws.Range("M4").Value = 80

# Row 6
$ws.Range("D6").Value = 44322
$ws.Range("M6").Value = 60

# Row 7
$ws.Range("D7").Value = 44327
$ws.Range("M7").Value = 60

# Row 8
$ws.Range("D8").Value = 44330

# Row 9
$ws.Range("D9").Value = 44316
$ws.Range("M9").Value = 120

# Row 10
$ws.Range("D10").Value = 44323
$ws.Range("M10").Value = 80
